$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update notes cell J8 -> "N/A"
$ws.Range("J8").Value = "N/A"

# Fill in row 11 with new data
$ws.Range("D11").NumberFormat = "mm-dd-yy"
$ws.Range("D11").Value = 45905
$ws.Range("E11").Value = 359
$ws.Range("F11").Value = 424
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 300
$ws.Range("J11").Value = "Este fin se vienen las fotos que no pude adelantar en la semana (rafael) "

# Update the selection
$ws.Range("H19").Select()
